$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 7).Value = 14.972416
$ws.Cells.Item(2, 8).Value = 29.944832
$ws.Cells.Item(2, 9).Value = 0.1033656722518705
$ws.Cells.Item(2, 10).Value = 0.083322905738039
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 13).Value = 1.631588
$ws.Cells.Item(2, 14).Value = 3.263176
$ws.Cells.Item(2, 15).Value = 0.06882101692447344
$ws.Cells.Item(2, 16).Value = 0.04759911448542056
$ws.Cells.Item(2, 17).Value = 24.428814276608
$ws.Cells.Item(2, 18).Value = 97.715257106432
$ws.Cells.Item(2, 19).Value = 0.007113730679455553
$ws.Cells.Item(2, 20).Value = 0.003966096529482824
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 7).Value = 14.972416
$ws.Cells.Item(3, 8).Value = 29.944832
$ws.Cells.Item(3, 9).Value = 0.1033656722518705
$ws.Cells.Item(3, 10).Value = 0.083322905738039
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.9053026666666666
$ws.Cells.Item(3, 14).Value = 2.715908
$ws.Cells.Item(3, 15).Value = 0.03818601886287323
$ws.Cells.Item(3, 16).Value = 0.03961625601066862
$ws.Cells.Item(3, 17).Value = 13.55456813124267
$ws.Cells.Item(3, 18).Value = 81.327408787456
$ws.Cells.Item(3, 19).Value = 0.003947123510383498
$ws.Cells.Item(3, 20).Value = 0.003300941565270962
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 7).Value = 14.972416
$ws.Cells.Item(4, 8).Value = 29.944832
$ws.Cells.Item(4, 9).Value = 0.1033656722518705
$ws.Cells.Item(4, 10).Value = 0.083322905738039
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 10.67936233333334
$ws.Cells.Item(4, 14).Value = 32.038087
$ws.Cells.Item(4, 15).Value = 0.450459660088771
$ws.Cells.Item(4, 16).Value = 0.4673313885021416
$ws.Cells.Item(4, 17).Value = 159.8958554693974
$ws.Cells.Item(4, 18).Value = 959.3751328163842
$ws.Cells.Item(4, 19).Value = 0.04656206558742488
$ws.Cells.Item(4, 20).Value = 0.03893940923259082
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 7).Value = 14.972416
$ws.Cells.Item(5, 8).Value = 29.944832
$ws.Cells.Item(5, 9).Value = 0.1033656722518705
$ws.Cells.Item(5, 10).Value = 0.083322905738039
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 8.307871
$ws.Cells.Item(5, 14).Value = 24.923613
$ws.Cells.Item(5, 15).Value = 0.3504292325619839
$ws.Cells.Item(5, 16).Value = 0.3635543741978111
$ws.Cells.Item(5, 17).Value = 124.388900686336
$ws.Cells.Item(5, 18).Value = 746.333404118016
$ws.Cells.Item(5, 19).Value = 0.03622235320047652
$ws.Cells.Item(5, 20).Value = 0.03029240685193597
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 7).Value = 14.972416
$ws.Cells.Item(6, 8).Value = 29.944832
$ws.Cells.Item(6, 9).Value = 0.1033656722518705
$ws.Cells.Item(6, 10).Value = 0.083322905738039
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 1.247457666666667
$ws.Cells.Item(6, 14).Value = 3.742373
$ws.Cells.Item(6, 15).Value = 0.05261824994436759
$ws.Cells.Item(6, 16).Value = 0.05458903867708846
$ws.Cells.Item(6, 17).Value = 18.67745512772267
$ws.Cells.Item(6, 18).Value = 112.064730766336
$ws.Cells.Item(6, 19).Value = 0.005438920778216502
$ws.Cells.Item(6, 20).Value = 0.004548517324021207
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 7).Value = 14.972416
$ws.Cells.Item(7, 8).Value = 29.944832
$ws.Cells.Item(7, 9).Value = 0.1033656722518705
$ws.Cells.Item(7, 10).Value = 0.083322905738039
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 13).Value = 0.936118
$ws.Cells.Item(7, 14).Value = 1.872236
$ws.Cells.Item(7, 15).Value = 0.03948582161753104
$ws.Cells.Item(7, 16).Value = 0.02730982812686961
$ws.Cells.Item(7, 17).Value = 14.015948121088
$ws.Cells.Item(7, 18).Value = 56.06379248435201
$ws.Cells.Item(7, 19).Value = 0.004081478495913535
$ws.Cells.Item(7, 20).Value = 0.002275534234737202
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 40.41312266666666
$ws.Cells.Item(8, 8).Value = 121.239368
$ws.Cells.Item(8, 9).Value = 0.2790017050179012
$ws.Cells.Item(8, 10).Value = 0.3373542530344942
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 13).Value = 1.631588
$ws.Cells.Item(8, 14).Value = 3.263176
$ws.Cells.Item(8, 15).Value = 0.06882101692447344
$ws.Cells.Item(8, 16).Value = 0.04759911448542056
$ws.Cells.Item(8, 17).Value = 65.93756598546132
$ws.Cells.Item(8, 18).Value = 395.6253959127679
$ws.Cells.Item(8, 19).Value = 0.01920118106299392
$ws.Cells.Item(8, 20).Value = 0.01605776371233242
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 40.41312266666666
$ws.Cells.Item(9, 8).Value = 121.239368
$ws.Cells.Item(9, 9).Value = 0.2790017050179012
$ws.Cells.Item(9, 10).Value = 0.3373542530344942
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.9053026666666666
$ws.Cells.Item(9, 14).Value = 2.715908
$ws.Cells.Item(9, 15).Value = 0.03818601886287323
$ws.Cells.Item(9, 16).Value = 0.03961625601066862
$ws.Cells.Item(9, 17).Value = 36.58610771846043
$ws.Cells.Item(9, 18).Value = 329.2749694661439
$ws.Cells.Item(9, 19).Value = 0.01065396437058737
$ws.Cells.Item(9, 20).Value = 0.0133647124545024
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 40.41312266666666
$ws.Cells.Item(10, 8).Value = 121.239368
$ws.Cells.Item(10, 9).Value = 0.2790017050179012
$ws.Cells.Item(10, 10).Value = 0.3373542530344942
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 10.67936233333334
$ws.Cells.Item(10, 14).Value = 32.038087
$ws.Cells.Item(10, 15).Value = 0.450459660088771
$ws.Cells.Item(10, 16).Value = 0.4673313885021416
$ws.Cells.Item(10, 17).Value = 431.5863799787796
$ws.Cells.Item(10, 18).Value = 3884.277419809016
$ws.Cells.Item(10, 19).Value = 0.1256790132065513
$ws.Cells.Item(10, 20).Value = 0.157656231487713
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 40.41312266666666
$ws.Cells.Item(11, 8).Value = 121.239368
$ws.Cells.Item(11, 9).Value = 0.2790017050179012
$ws.Cells.Item(11, 10).Value = 0.3373542530344942
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 8.307871
$ws.Cells.Item(11, 14).Value = 24.923613
$ws.Cells.Item(11, 15).Value = 0.3504292325619839
$ws.Cells.Item(11, 16).Value = 0.3635543741978111
$ws.Cells.Item(11, 17).Value = 335.7470098218426
$ws.Cells.Item(11, 18).Value = 3021.723088396584
$ws.Cells.Item(11, 19).Value = 0.09777035337290811
$ws.Cells.Item(11, 20).Value = 0.1226466143449255
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 40.41312266666666
$ws.Cells.Item(12, 8).Value = 121.239368
$ws.Cells.Item(12, 9).Value = 0.2790017050179012
$ws.Cells.Item(12, 10).Value = 0.3373542530344942
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 1.247457666666667
$ws.Cells.Item(12, 14).Value = 3.742373
$ws.Cells.Item(12, 15).Value = 0.05261824994436759
$ws.Cells.Item(12, 16).Value = 0.05458903867708846
$ws.Cells.Item(12, 17).Value = 50.41365970447377
$ws.Cells.Item(12, 18).Value = 453.722937340264
$ws.Cells.Item(12, 19).Value = 0.01468058144953664
$ws.Cells.Item(12, 20).Value = 0.01841584436678029
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 40.41312266666666
$ws.Cells.Item(13, 8).Value = 121.239368
$ws.Cells.Item(13, 9).Value = 0.2790017050179012
$ws.Cells.Item(13, 10).Value = 0.3373542530344942
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 13).Value = 0.936118
$ws.Cells.Item(13, 14).Value = 1.872236
$ws.Cells.Item(13, 15).Value = 0.03948582161753104
$ws.Cells.Item(13, 16).Value = 0.02730982812686961
$ws.Cells.Item(13, 17).Value = 37.83145156447466
$ws.Cells.Item(13, 18).Value = 226.988709386848
$ws.Cells.Item(13, 19).Value = 0.01101661155532386
$ws.Cells.Item(13, 20).Value = 0.009213086668240516
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 5.007042333333334
$ws.Cells.Item(14, 8).Value = 15.021127
$ws.Cells.Item(14, 9).Value = 0.03456732011577652
$ws.Cells.Item(14, 10).Value = 0.04179699352128983
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 13).Value = 1.631588
$ws.Cells.Item(14, 14).Value = 3.263176
$ws.Cells.Item(14, 15).Value = 0.06882101692447344
$ws.Cells.Item(14, 16).Value = 0.04759911448542056
$ws.Cells.Item(14, 17).Value = 8.169430186558667
$ws.Cells.Item(14, 18).Value = 49.016581119352
$ws.Cells.Item(14, 19).Value = 0.002378958122721547
$ws.Cells.Item(14, 20).Value = 0.001989499879766256
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 5.007042333333334
$ws.Cells.Item(15, 8).Value = 15.021127
$ws.Cells.Item(15, 9).Value = 0.03456732011577652
$ws.Cells.Item(15, 10).Value = 0.04179699352128983
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 0.9053026666666666
$ws.Cells.Item(15, 14).Value = 2.715908
$ws.Cells.Item(15, 15).Value = 0.03818601886287323
$ws.Cells.Item(15, 16).Value = 0.03961625601066862
$ws.Cells.Item(15, 17).Value = 4.532888776479555
$ws.Cells.Item(15, 18).Value = 40.79599898831599
$ws.Cells.Item(15, 19).Value = 0.001319988337980019
$ws.Cells.Item(15, 20).Value = 0.001655840395815676
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 5.007042333333334
$ws.Cells.Item(16, 8).Value = 15.021127
$ws.Cells.Item(16, 9).Value = 0.03456732011577652
$ws.Cells.Item(16, 10).Value = 0.04179699352128983
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 10.67936233333334
$ws.Cells.Item(16, 14).Value = 32.038087
$ws.Cells.Item(16, 15).Value = 0.450459660088771
$ws.Cells.Item(16, 16).Value = 0.4673313885021416
$ws.Cells.Item(16, 17).Value = 53.47201929600546
$ws.Cells.Item(16, 18).Value = 481.2481736640491
$ws.Cells.Item(16, 19).Value = 0.01557118326953243
$ws.Cells.Item(16, 20).Value = 0.01953304701751939
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 5.007042333333334
$ws.Cells.Item(17, 8).Value = 15.021127
$ws.Cells.Item(17, 9).Value = 0.03456732011577652
$ws.Cells.Item(17, 10).Value = 0.04179699352128983
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 8.307871
$ws.Cells.Item(17, 14).Value = 24.923613
$ws.Cells.Item(17, 15).Value = 0.3504292325619839
$ws.Cells.Item(17, 16).Value = 0.3635543741978111
$ws.Cells.Item(17, 17).Value = 41.59786179687234
$ws.Cells.Item(17, 18).Value = 374.380756171851
$ws.Cells.Item(17, 19).Value = 0.01211339945989599
$ws.Cells.Item(17, 20).Value = 0.01519547982298249
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 5.007042333333334
$ws.Cells.Item(18, 8).Value = 15.021127
$ws.Cells.Item(18, 9).Value = 0.03456732011577652
$ws.Cells.Item(18, 10).Value = 0.04179699352128983
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 1.247457666666667
$ws.Cells.Item(18, 14).Value = 3.742373
$ws.Cells.Item(18, 15).Value = 0.05261824994436759
$ws.Cells.Item(18, 16).Value = 0.05458903867708846
$ws.Cells.Item(18, 17).Value = 6.246073346041223
$ws.Cells.Item(18, 18).Value = 56.214660114371
$ws.Cells.Item(18, 19).Value = 0.001818871889758894
$ws.Cells.Item(18, 20).Value = 0.002281657695919706
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 5.007042333333334
$ws.Cells.Item(19, 8).Value = 15.021127
$ws.Cells.Item(19, 9).Value = 0.03456732011577652
$ws.Cells.Item(19, 10).Value = 0.04179699352128983
$ws.Cells.Item(19, 11).Value = 2
$ws.Cells.Item(19, 13).Value = 0.936118
$ws.Cells.Item(19, 14).Value = 1.872236
$ws.Cells.Item(19, 15).Value = 0.03948582161753104
$ws.Cells.Item(19, 16).Value = 0.02730982812686961
$ws.Cells.Item(19, 17).Value = 4.687182454995334
$ws.Cells.Item(19, 18).Value = 28.123094729972
$ws.Cells.Item(19, 19).Value = 0.001364919035887644
$ws.Cells.Item(19, 20).Value = 0.001141468709286308
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 12.521722
$ws.Cells.Item(20, 8).Value = 37.565166
$ws.Cells.Item(20, 9).Value = 0.08644671723528362
$ws.Cells.Item(20, 10).Value = 0.1045268440862112
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 13).Value = 1.631588
$ws.Cells.Item(20, 14).Value = 3.263176
$ws.Cells.Item(20, 15).Value = 0.06882101692447344
$ws.Cells.Item(20, 16).Value = 0.04759911448542056
$ws.Cells.Item(20, 17).Value = 20.430291354536
$ws.Cells.Item(20, 18).Value = 122.581748127216
$ws.Cells.Item(20, 19).Value = 0.005949350989914624
$ws.Cells.Item(20, 20).Value = 0.004975385218459271
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 12.521722
$ws.Cells.Item(21, 8).Value = 37.565166
$ws.Cells.Item(21, 9).Value = 0.08644671723528362
$ws.Cells.Item(21, 10).Value = 0.1045268440862112
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 13).Value = 0.9053026666666666
$ws.Cells.Item(21, 14).Value = 2.715908
$ws.Cells.Item(21, 15).Value = 0.03818601886287323
$ws.Cells.Item(21, 16).Value = 0.03961625601066862
$ws.Cells.Item(21, 17).Value = 11.33594831785867
$ws.Cells.Item(21, 18).Value = 102.023534860728
$ws.Cells.Item(21, 19).Value = 0.003301055974980008
$ws.Cells.Item(21, 20).Value = 0.004140962215306586
$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 7).Value = 12.521722
$ws.Cells.Item(22, 8).Value = 37.565166
$ws.Cells.Item(22, 9).Value = 0.08644671723528362
$ws.Cells.Item(22, 10).Value = 0.1045268440862112
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 10.67936233333334
$ws.Cells.Item(22, 14).Value = 32.038087
$ws.Cells.Item(22, 15).Value = 0.450459660088771
$ws.Cells.Item(22, 16).Value = 0.4673313885021416
$ws.Cells.Item(22, 17).Value = 133.7240062752714
$ws.Cells.Item(22, 18).Value = 1203.516056477442
$ws.Cells.Item(22, 19).Value = 0.03894075886159596
$ws.Cells.Item(22, 20).Value = 0.04884867518255594
$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 7).Value = 12.521722
$ws.Cells.Item(23, 8).Value = 37.565166
$ws.Cells.Item(23, 9).Value = 0.08644671723528362
$ws.Cells.Item(23, 10).Value = 0.1045268440862112
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 8.307871
$ws.Cells.Item(23, 14).Value = 24.923613
$ws.Cells.Item(23, 15).Value = 0.3504292325619839
$ws.Cells.Item(23, 16).Value = 0.3635543741978111
$ws.Cells.Item(23, 17).Value = 104.028851073862
$ws.Cells.Item(23, 18).Value = 936.2596596647581
$ws.Cells.Item(23, 19).Value = 0.03029345677826326
$ws.Cells.Item(23, 20).Value = 0.03800119138863468
$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 7).Value = 12.521722
$ws.Cells.Item(24, 8).Value = 37.565166
$ws.Cells.Item(24, 9).Value = 0.08644671723528362
$ws.Cells.Item(24, 10).Value = 0.1045268440862112
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 1.247457666666667
$ws.Cells.Item(24, 14).Value = 3.742373
$ws.Cells.Item(24, 15).Value = 0.05261824994436759
$ws.Cells.Item(24, 16).Value = 0.05458903867708846
$ws.Cells.Item(24, 17).Value = 15.62031810876867
$ws.Cells.Item(24, 18).Value = 140.582862978918
$ws.Cells.Item(24, 19).Value = 0.004548674974356223
$ws.Cells.Item(24, 20).Value = 0.005706019934616178
$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 7).Value = 12.521722
$ws.Cells.Item(25, 8).Value = 37.565166
$ws.Cells.Item(25, 9).Value = 0.08644671723528362
$ws.Cells.Item(25, 10).Value = 0.1045268440862112
$ws.Cells.Item(25, 11).Value = 2
$ws.Cells.Item(25, 13).Value = 0.936118
$ws.Cells.Item(25, 14).Value = 1.872236
$ws.Cells.Item(25, 15).Value = 0.03948582161753104
$ws.Cells.Item(25, 16).Value = 0.02730982812686961
$ws.Cells.Item(25, 17).Value = 11.721809355196
$ws.Cells.Item(25, 18).Value = 70.330856131176
$ws.Cells.Item(25, 19).Value = 0.003413419656173555
$ws.Cells.Item(25, 20).Value = 0.002854610146638524
$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 7).Value = 11.74303333333333
$ws.Cells.Item(26, 8).Value = 35.2291
$ws.Cells.Item(26, 9).Value = 0.08107085287879548
$ws.Cells.Item(26, 10).Value = 0.09802663038937569
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 13).Value = 1.631588
$ws.Cells.Item(26, 14).Value = 3.263176
$ws.Cells.Item(26, 15).Value = 0.06882101692447344
$ws.Cells.Item(26, 16).Value = 0.04759911448542056
$ws.Cells.Item(26, 17).Value = 19.15979227026667
$ws.Cells.Item(26, 18).Value = 114.9587536216
$ws.Cells.Item(26, 19).Value = 0.00557937853805308
$ws.Cells.Item(26, 20).Value = 0.0046659808025239
$ws.Cells.Item(27, 5).Value = 3
$ws.Cells.Item(27, 7).Value = 11.74303333333333
$ws.Cells.Item(27, 8).Value = 35.2291
$ws.Cells.Item(27, 9).Value = 0.08107085287879548
$ws.Cells.Item(27, 10).Value = 0.09802663038937569
$ws.Cells.Item(27, 11).Value = 3
$ws.Cells.Item(27, 13).Value = 0.9053026666666666
$ws.Cells.Item(27, 14).Value = 2.715908
$ws.Cells.Item(27, 15).Value = 0.03818601886287323
$ws.Cells.Item(27, 16).Value = 0.03961625601066862
$ws.Cells.Item(27, 17).Value = 10.63099939142222
$ws.Cells.Item(27, 18).Value = 95.6789945228
$ws.Cells.Item(27, 19).Value = 0.003095773117258904
$ws.Cells.Item(27, 20).Value = 0.003883448085368696
$ws.Cells.Item(28, 5).Value = 3
$ws.Cells.Item(28, 7).Value = 11.74303333333333
$ws.Cells.Item(28, 8).Value = 35.2291
$ws.Cells.Item(28, 9).Value = 0.08107085287879548
$ws.Cells.Item(28, 10).Value = 0.09802663038937569
$ws.Cells.Item(28, 11).Value = 3
$ws.Cells.Item(28, 13).Value = 10.67936233333334
$ws.Cells.Item(28, 14).Value = 32.038087
$ws.Cells.Item(28, 15).Value = 0.450459660088771
$ws.Cells.Item(28, 16).Value = 0.4673313885021416
$ws.Cells.Item(28, 17).Value = 125.4081078590778
$ws.Cells.Item(28, 18).Value = 1128.6729707317
$ws.Cells.Item(28, 19).Value = 0.03651914883088897
$ws.Cells.Item(28, 20).Value = 0.04581092129005317
$ws.Cells.Item(29, 5).Value = 3
$ws.Cells.Item(29, 7).Value = 11.74303333333333
$ws.Cells.Item(29, 8).Value = 35.2291
$ws.Cells.Item(29, 9).Value = 0.08107085287879548
$ws.Cells.Item(29, 10).Value = 0.09802663038937569
$ws.Cells.Item(29, 11).Value = 3
$ws.Cells.Item(29, 13).Value = 8.307871
$ws.Cells.Item(29, 14).Value = 24.923613
$ws.Cells.Item(29, 15).Value = 0.3504292325619839
$ws.Cells.Item(29, 16).Value = 0.3635543741978111
$ws.Cells.Item(29, 17).Value = 97.55960608203335
$ws.Cells.Item(29, 18).Value = 878.0364547383
$ws.Cells.Item(29, 19).Value = 0.0284095967574618
$ws.Cells.Item(29, 20).Value = 0.03563801026592961
$ws.Cells.Item(30, 5).Value = 3
$ws.Cells.Item(30, 7).Value = 11.74303333333333
$ws.Cells.Item(30, 8).Value = 35.2291
$ws.Cells.Item(30, 9).Value = 0.08107085287879548
$ws.Cells.Item(30, 10).Value = 0.09802663038937569
$ws.Cells.Item(30, 11).Value = 3
$ws.Cells.Item(30, 13).Value = 1.247457666666667
$ws.Cells.Item(30, 14).Value = 3.742373
$ws.Cells.Item(30, 15).Value = 0.05261824994436759
$ws.Cells.Item(30, 16).Value = 0.05458903867708846
$ws.Cells.Item(30, 17).Value = 14.64893696158889
$ws.Cells.Item(30, 18).Value = 131.8404326543
$ws.Cells.Item(30, 19).Value = 0.004265806399979514
$ws.Cells.Item(30, 20).Value = 0.005351179517710285
$ws.Cells.Item(31, 5).Value = 3
$ws.Cells.Item(31, 7).Value = 11.74303333333333
$ws.Cells.Item(31, 8).Value = 35.2291
$ws.Cells.Item(31, 9).Value = 0.08107085287879548
$ws.Cells.Item(31, 10).Value = 0.09802663038937569
$ws.Cells.Item(31, 11).Value = 2
$ws.Cells.Item(31, 13).Value = 0.936118
$ws.Cells.Item(31, 14).Value = 1.872236
$ws.Cells.Item(31, 15).Value = 0.03948582161753104
$ws.Cells.Item(31, 16).Value = 0.02730982812686961
$ws.Cells.Item(31, 17).Value = 10.99286487793333
$ws.Cells.Item(31, 18).Value = 65.9571892676
$ws.Cells.Item(31, 19).Value = 0.003201149235153221
$ws.Cells.Item(31, 20).Value = 0.002677090427790023
$ws.Cells.Item(32, 5).Value = 2
$ws.Cells.Item(32, 7).Value = 60.1916805
$ws.Cells.Item(32, 8).Value = 120.383361
$ws.Cells.Item(32, 9).Value = 0.4155477325003729
$ws.Cells.Item(32, 10).Value = 0.3349723732305901
$ws.Cells.Item(32, 11).Value = 2
$ws.Cells.Item(32, 13).Value = 1.631588
$ws.Cells.Item(32, 14).Value = 3.263176
$ws.Cells.Item(32, 15).Value = 0.06882101692447344
$ws.Cells.Item(32, 16).Value = 0.04759911448542056
$ws.Cells.Item(32, 17).Value = 98.208023603634
$ws.Cells.Item(32, 18).Value = 392.832094414536
$ws.Cells.Item(32, 19).Value = 0.02859841753133472
$ws.Cells.Item(32, 20).Value = 0.01594438834285588
$ws.Cells.Item(33, 5).Value = 2
$ws.Cells.Item(33, 7).Value = 60.1916805
$ws.Cells.Item(33, 8).Value = 120.383361
$ws.Cells.Item(33, 9).Value = 0.4155477325003729
$ws.Cells.Item(33, 10).Value = 0.3349723732305901
$ws.Cells.Item(33, 11).Value = 3
$ws.Cells.Item(33, 13).Value = 0.9053026666666666
$ws.Cells.Item(33, 14).Value = 2.715908
$ws.Cells.Item(33, 15).Value = 0.03818601886287323
$ws.Cells.Item(33, 16).Value = 0.03961625601066862
$ws.Cells.Item(33, 17).Value = 54.49168886779799
$ws.Cells.Item(33, 18).Value = 326.9501332067879
$ws.Cells.Item(33, 19).Value = 0.01586811355168343
$ws.Cells.Item(33, 20).Value = 0.0132703512944043
$ws.Cells.Item(34, 5).Value = 2
$ws.Cells.Item(34, 7).Value = 60.1916805
$ws.Cells.Item(34, 8).Value = 120.383361
$ws.Cells.Item(34, 9).Value = 0.4155477325003729
$ws.Cells.Item(34, 10).Value = 0.3349723732305901
$ws.Cells.Item(34, 11).Value = 3
$ws.Cells.Item(34, 13).Value = 10.67936233333334
$ws.Cells.Item(34, 14).Value = 32.038087
$ws.Cells.Item(34, 15).Value = 0.450459660088771
$ws.Cells.Item(34, 16).Value = 0.4673313885021416
$ws.Cells.Item(34, 17).Value = 642.8087655117346
$ws.Cells.Item(34, 18).Value = 3856.852593070407
$ws.Cells.Item(34, 19).Value = 0.1871874903327775
$ws.Cells.Item(34, 20).Value = 0.1565431042917093
$ws.Cells.Item(35, 5).Value = 2
$ws.Cells.Item(35, 7).Value = 60.1916805
$ws.Cells.Item(35, 8).Value = 120.383361
$ws.Cells.Item(35, 9).Value = 0.4155477325003729
$ws.Cells.Item(35, 10).Value = 0.3349723732305901
$ws.Cells.Item(35, 11).Value = 3
$ws.Cells.Item(35, 13).Value = 8.307871
$ws.Cells.Item(35, 14).Value = 24.923613
$ws.Cells.Item(35, 15).Value = 0.3504292325619839
$ws.Cells.Item(35, 16).Value = 0.3635543741978111
$ws.Cells.Item(35, 17).Value = 500.0647168672155
$ws.Cells.Item(35, 18).Value = 3000.388301203293
$ws.Cells.Item(35, 19).Value = 0.1456200729929782
$ws.Cells.Item(35, 20).Value = 0.1217806715234028
$ws.Cells.Item(36, 5).Value = 2
$ws.Cells.Item(36, 7).Value = 60.1916805
$ws.Cells.Item(36, 8).Value = 120.383361
$ws.Cells.Item(36, 9).Value = 0.4155477325003729
$ws.Cells.Item(36, 10).Value = 0.3349723732305901
$ws.Cells.Item(36, 11).Value = 3
$ws.Cells.Item(36, 13).Value = 1.247457666666667
$ws.Cells.Item(36, 14).Value = 3.742373
$ws.Cells.Item(36, 15).Value = 0.05261824994436759
$ws.Cells.Item(36, 16).Value = 0.05458903867708846
$ws.Cells.Item(36, 17).Value = 75.08657330927551
$ws.Cells.Item(36, 18).Value = 450.519439855653
$ws.Cells.Item(36, 19).Value = 0.02186539445251982
$ws.Cells.Item(36, 20).Value = 0.01828581983804079
$ws.Cells.Item(37, 5).Value = 2
$ws.Cells.Item(37, 7).Value = 60.1916805
$ws.Cells.Item(37, 8).Value = 120.383361
$ws.Cells.Item(37, 9).Value = 0.4155477325003729
$ws.Cells.Item(37, 10).Value = 0.3349723732305901
$ws.Cells.Item(37, 11).Value = 2
$ws.Cells.Item(37, 13).Value = 0.936118
$ws.Cells.Item(37, 14).Value = 1.872236
$ws.Cells.Item(37, 15).Value = 0.03948582161753104
$ws.Cells.Item(37, 16).Value = 0.02730982812686961
$ws.Cells.Item(37, 17).Value = 56.346515566299
$ws.Cells.Item(37, 18).Value = 225.386062265196
$ws.Cells.Item(37, 19).Value = 0.01640824363907923
$ws.Cells.Item(37, 20).Value = 0.009148037940177033
